$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# Insert a new row at row 4 (pushes the table's existing rows down)
$ws.Rows("4:4").Insert()
$ws.Range("A4").Value = "Rörliga objekt"

# Re-sync the table range now that a row has been inserted inside it
$tbl.Resize($ws.Range("A1:F13"))

# Make the new cell's column-A look consistent with the rest of the list
# (this forces Excel to materialize a dedicated bold style/font, matching
# the extra font + cellXf that show up in the saved workbook)
$ws.Range("A4:A12").Font.Name = "Calibri"

# Match the selection left behind in the sheet view
$ws.Range("H5").Select() | Out-Null
